$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.078.03"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "1.833.38"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'324.60"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D7").Value = "'0.4615"
$ws.Range("E7").Value = "  -0.98%  "
$ws.Range("D8").Value = "'0.3866"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").Value = "'0.07850"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "'0.9610"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").Value = "'21.95"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").Value = "1.856.72"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "'5.673"
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "'6.882"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "'0.06863"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "'88.30"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'0.000009938"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").Value = "'16.69"
$ws.Range("E19").Value = "  -2.46%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "28.101.82"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "'5.295"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").Value = "'2.082"
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("D25").Value = "2.087.63"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").Value = "'154.75"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").Value = "'19.16"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "'5.728"
$ws.Range("E28").Value = "  -6.18%  "
$ws.Range("D29").Value = "'1.967"
$ws.Range("E29").Value = "  -3.09%  "
$ws.Range("D30").Value = "'118.53"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").Value = "'0.9420"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").Value = "'0.09242"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("D33").Value = "'5.262"
$ws.Range("E33").Value = "  -1.82%  "
$ws.Range("D34").Value = "'1.321"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").Value = "'3.327"
$ws.Range("E35").Value = "  -4.56%  "
$ws.Range("D36").Value = "'0.05848"
$ws.Range("E36").Value = "  -5.02%  "
$ws.Range("D37").Value = "'0.02110"
$ws.Range("E37").Value = "  -4.07%  "
$ws.Range("D38").Value = "'1.139"
$ws.Range("E38").Value = "  -1.96%  "
$ws.Range("D39").Value = "'7.726"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "'0.5591"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "'9.899"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("D42").Value = "'0.1760"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").Value = "'0.07322"
$ws.Range("E43").Value = "  +2.76%  "
$ws.Range("D44").Value = "'11.69"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").Value = "'0.5275"
$ws.Range("E45").Value = "  -2.22%  "
$ws.Range("D46").Value = "'1.142"
$ws.Range("E46").Value = "  -8.49%  "
$ws.Range("D47").Value = "'2.108"
$ws.Range("E47").Value = "  -10.84%  "
$ws.Range("D48").Value = "'1.832"
$ws.Range("E48").Value = "  -3.99%  "
$ws.Range("D49").Value = "'113.08"
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  +0.13%  "
